$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B, C, D are treated as plain text so numeric-looking
# values (e.g. "1.001", "0.07370") are preserved exactly as strings,
# matching the original inline-string cell content in the workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.898.53"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.815.82"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.72"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4662"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3697"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07370"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8708"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.850.02"
$ws.Range("E12").Value = "  +2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.374"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.35"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07087"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.502"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008711"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.925.71"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.339"
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.040.06"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.891"
$ws.Range("E25").Value = "  -1.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.76"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.197"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.36"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.311"
$ws.Range("E29").Value = "  +0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.74"
$ws.Range("E30").Value = "  -0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08918"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7657"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.481"
$ws.Range("E34").Value = "  -0.89%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.090"
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01961"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05280"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.961"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5350"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.237"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.365"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1663"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.440"
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4935"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.40"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.671"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.76"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06284"
$ws.Range("E51").Value = "  -0.73%  "
